$d = $word.ActiveDocument

# 1. Split the "Trabalho feito ... aulas." sentence and insert the
#    author's name in the middle, keeping the text identical aside
#    from the inserted segment.
$d.Content.Find.Execute(
    "Trabalho feito para a matéria",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Trabalho feito pelo aluno Pedro Graça Carneiro para a matéria",
    2)

# 2. Remove the old "_GoBack" bookmark (it used to sit right after
#    " ID).") ...
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ... and re-create it right after the newly inserted author name,
# i.e. right before " para a matéria de POO em Java ...".
$rng = $d.Content
$rng.Find.Execute("Pedro Graça Carneiro", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkRange = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null
